$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.110.59'
$ws.Range("E2").Value = '  +1.54%  '

$ws.Range("D3").Value = '3.601.50'
$ws.Range("E3").Value = '  +0.28%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '659.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.72%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.68'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +15.19%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.415'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.42%  '

$ws.Range("E9").Value = '  +5.78%  '

$ws.Range("E10").Value = '  +0.03%  '

$ws.Range("D11").Value = '3.600.04'
$ws.Range("E11").Value = '  +0.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.80'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.38%  '

$ws.Range("E13").Value = '  +1.36%  '

$ws.Range("E14").Value = '  +0.12%  '

$ws.Range("D15").Value = '4.271.24'
$ws.Range("E15").Value = '  +0.18%  '

$ws.Range("D16").Value = '96.715.61'
$ws.Range("E16").Value = '  +1.35%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000259'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.46%  '

$ws.Range("D18").Value = '3.596.18'
$ws.Range("E18").Value = '  +0.16%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.25%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.92%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.536'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +12.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '512.30'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.43'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.38%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000204'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.33%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '98.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.15%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.51%  '

$ws.Range("D29").Value = '3.793.24'
$ws.Range("E29").Value = '  +0.19%  '

$ws.Range("E30").Value = '  -0.61%  '

$ws.Range("E31").Value = '  +8.44%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.66'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.82%  '

$ws.Range("E33").Value = '  +0.09%  '

$ws.Range("E34").Value = '  +5.29%  '

$ws.Range("E35").Value = '  +0.28%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.77'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.42%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '623.34'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.86%  '

$ws.Range("E38").Value = '  +2.02%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.74'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.81%  '

$ws.Range("E40").Value = '  +11.14%  '

$ws.Range("E41").Value = '  +2.10%  '

$ws.Range("B42").Value = 'ImmutableX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.93'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +10.30%  '

$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.922'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.93'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0433'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '33.32'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.75%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.46%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.53'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.74%  '
